# Generate Report for Handoff
# Replace the stale run's file id / xliff hashes / timestamps with the
# values produced by the latest handoff run.

$wb = $excel.ActiveWorkbook

$newId   = "defb6900-bdec-4f8e-9578-e137b318b59a"
$newHash = "014b33f13e197f0de5091d681689e2f60feaf7a3"

$newHoDate = "2016-08-18 12:58:13"
$newZhDate = "2016-08-18 12:58:01"

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = "$newId.md"
$wsOverview.Range("G2").Value2 = $newHoDate
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newId.md"
}
$wsOverview.Range("B2").Value2 = "e2e\$newId.md"

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("G2").Value2 = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value2 = $newZhDate
foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
$wsZh.Range("A2").Value2 = "$newId.md"

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value2 = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value2 = $newHoDate
foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
$wsDe.Range("A2").Value2 = "$newId.md"

Write-Host "Updated localization status report for handoff id $newId"
